# Adds two new "demography and housing model" parameter rows
# (prop_less_ims / prop_more_ims) to the parameter sheet, mirroring the
# existing block-style rows (e.g. rows 180/181) for look & feel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Stamp formatting for the two new rows by copying it from an
#        existing, similarly-styled block (rows 180:181 use the same
#        category/value/unit/relevance/description layout we need) ---
$ws.Range("A180:J181").Copy()
$ws.Range("A185:J186").PasteSpecial(-4122)

# Column C needs the wrapped "description" style (like column G) rather
# than the plain "no parameter" style it would inherit from C180:C181.
$ws.Range("G180:G181").Copy()
$ws.Range("C185:C186").PasteSpecial(-4122)

# Match the taller row height used for these longer descriptions.
$ws.Range("A185:J186").RowHeight = 51

# --- 2. Fill in the new cell text in the same order the original
#        authoring tool used, so brand-new shared-string entries land in
#        the same index order as the canonical edit (A185, B185, G185,
#        G186, C185, B186, then the remaining reused/number cells). ---
$ws.Range("A185").Value = "demography and housing model"
$ws.Range("B185").Value = "prop_less_ims"
$ws.Range("G185").Value = "less living space available (according to reserves, compared to trend calcuation). Then immigration* is descreased, and emigration* increased. This parameter determines how much of the difference is compensated by changes of immigration* (the remainder is corrected with emigration*)"
$ws.Range("G186").Value = "more living space available (according to reserves, compared to trend calcuation). Then immigration* is inscreased, and emigration* decreased. This parameter determines how much of the difference is compensated by changes of immigration* (the remainder is corrected with emigration*)"
$ws.Range("C185").Value = "similar to mod.ant.zuz (in the previous model only one parameter for both situations: less or more living space available)"
$ws.Range("B186").Value = "prop_more_ims"

# --- 3. Remaining cells: reused shared strings (category label, unit,
#        relevance) and plain numbers. ---
$ws.Range("A186").Value = "demography and housing model"
$ws.Range("C186").Value = "similar to mod.ant.zuz (in the previous model only one parameter for both situations: less or more living space available)"

$ws.Range("D185").Value = 80
$ws.Range("E185").Value = "percent"
$ws.Range("F185").Value = "no parameter"
$ws.Range("H185").Value = 80
$ws.Range("I185").Value = 80
$ws.Range("J185").Value = 80

$ws.Range("D186").Value = 60
$ws.Range("E186").Value = "percent"
$ws.Range("F186").Value = "no parameter"
$ws.Range("H186").Value = 60
$ws.Range("I186").Value = 60
$ws.Range("J186").Value = 60

# --- 4. Restore the top-row freeze (header only) and leave the
#        selection on the newly-edited area, similar to the author's
#        end state. ---
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A178").Select()
